# Fix typo'd requirement IDs (rg_ -> rq_) on the requirements sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("requirements")

$ws.Range("A80").Value = "rq_table_search"
$ws.Range("A81").Value = "rq_single_view"

# New requirements to append at the bottom of the requirements sheet
$newReqs = @(
    @("rq_spreadsheet_export", "VRM2 shall offer the option to export the selected specobjects to a spreadsheet file in xlsx format. ", 13),
    @("rq_spreadsheet_export_cfg", "VRM2 shall offer user selection of specobjects fields and their order in exported file.", 13),
    @("rq_vql_parents", "VRM2 shall have a VQL function that finds the set of parent specobjects from a selection set, and aplies a specified filter.", 13),
    @("rq_vql_children", "VRM2 shall have a VQL function that finds the set of children specobjects from a selection set, and aplies a specified filter.", 13),
    @("rq_vql_descendants", "VRM2 shall have a VQL function that finds the set of descendant specobjects from a selection set, and aplies a specified filter.", 13),
    @("rq_vql_ancestors", "VRM2 shall have a VQL function that finds the set of ancestor specobjects from a selection set, and aplies a specified filter.", 13),
    @("rq_doctype_filetypes", "VRM2 shall show the file types (from the <sourcefile>) that contribute to each doctype in the hierarchy view. I.e. a list of filetypes shall be added to each node in diagram.", 24)
)

$row = 82
foreach ($req in $newReqs) {
    $ws.Range("A$row").Value = $req[0]
    $ws.Range("B$row").Value = 1
    $ws.Range("C$row").Value = "approved"
    $ws.Range("D$row").Value = $req[1]
    $ws.Range("E$row").Value = "sourcecode;testcode"
    $ws.Rows.Item($row).RowHeight = $req[2]
    $row = $row + 1
}

# Drop the two trailing placeholder rows left over at the bottom of the sheet
$ws.Rows.Item(1048576).Delete() | Out-Null
$ws.Rows.Item(1048575).Delete() | Out-Null

# Keep selection / scroll position roughly in sync with the final document
$ws.Application.ActiveWindow.ScrollRow = 40
$ws.Range("D82").Select() | Out-Null

Write-Host "applied"
